$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.296.59'
$ws.Range('E2').Value = '  -4.82%  '
$ws.Range('D3').Value = '3.266.98'
$ws.Range('E3').Value = '  -6.97%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -11.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '3.259.10'
$ws.Range('E8').Value = '  -7.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.542'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -10.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -13.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.71'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.509'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -11.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -15.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000245'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -10.25%  '
$ws.Range('D15').Value = '3.785.67'
$ws.Range('E15').Value = '  -7.21%  '
$ws.Range('D16').Value = '67.254.55'
$ws.Range('E16').Value = '  -5.00%  '
$ws.Range('D17').Value = '3.265.21'
$ws.Range('E17').Value = '  -7.18%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '534.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -10.54%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.114'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -13.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -13.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.762'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -12.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.86'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -12.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -11.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -12.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -11.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -14.24%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '29.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.57%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -9.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.62'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -17.80%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.78'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -13.80%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '531.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -11.50%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0456'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.24'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0860'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -12.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -16.23%  '
$ws.Range('E41').Value = '  -11.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.78'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -16.69%  '
$ws.Range('D43').Value = '2.933.06'
$ws.Range('E43').Value = '  -11.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.268'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -12.35%  '
$ws.Range('D45').Value = '0.0₃0592'
$ws.Range('E45').Value = '  -17.16%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.95%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -14.69%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.35'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -16.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.37'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.25%  '
$ws.Range('E51').Value = '  -12.09%  '
